# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet (fund-holdings detail, same layout as
# the other quarterly sheets) right before the "总计" (totals) summary
# sheet, and prepends a corresponding "2022-Q1" row to the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q1" sheet, positioned right before "总计".
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totals)
$q1.Name = "2022-Q1"

# NOTE: `$totals` was resolved to a position-based handle, and inserting
# a sheet in front of it shifts what that handle now points at (it now
# refers to the freshly-inserted "2022-Q1" sheet). Re-resolve "总计" by
# name so subsequent writes land on the right sheet.
$totals = $wb.Worksheets.Item("总计")

# Use the existing "2021-Q4" sheet as a formatting template: its header
# row (B1:H1) and index column (A2:A?) both use the bold/bordered
# "header" style already present in the workbook.
$template = $wb.Worksheets.Item("2021-Q4")

$template.Range("B1:H1").Copy()
$q1.Range("B1").PasteSpecial(-4122)

$template.Range("A2").Copy()
$q1.Range("A2:A7").PasteSpecial(-4122)

# ---- headers ----------------------------------------------------------
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# ---- data rows ----------------------------------------------------------
# index, code, name, scale, stock position, position %, held value, rank
$q1Rows = @(
    @(0, "012671", "嘉实核心蓝筹混合型证券投资基金A", "11.11", "93.27", "5.72", "0.6355", 7),
    @(1, "009126", "嘉实基础产业优选股票A",           "2.17",  "90.79", "5.35", "0.1161", 5),
    @(2, "010783", "德邦沪港深龙头混合A",             "0.93",  "81.58", "3.81", "0.0354", 6),
    @(3, "012672", "嘉实核心蓝筹混合型证券投资基金C", "0.47",  "93.27", "5.72", "0.0269", 7),
    @(4, "010784", "德邦沪港深龙头混合C",             "0.27",  "81.58", "3.81", "0.0103", 6),
    @(5, "009127", "嘉实基础产业优选股票C",           "0.11",  "90.79", "5.35", "0.0059", 5)
)

$r = 2
foreach ($row in $q1Rows) {
    $q1.Range("A$r").Value = $row[0]
    $q1.Range("B$r").Value = "'" + $row[1]
    $q1.Range("C$r").Value = $row[2]
    $q1.Range("D$r").Value = "'" + $row[3]
    $q1.Range("E$r").Value = "'" + $row[4]
    $q1.Range("F$r").Value = "'" + $row[5]
    $q1.Range("G$r").Value = "'" + $row[6]
    $q1.Range("H$r").Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2) Prepend a "2022-Q1" row to the "总计" sheet, pushing the rest down.
# ---------------------------------------------------------------------
$totals.Rows("2:2").Insert()
$totals.Range("A2:D2").ClearFormats()

$totals.Range("A3").Copy()
$totals.Range("A2").PasteSpecial(-4122)

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q1"
$totals.Range("C2").Value = 6
$totals.Range("D2").Value = 0.83

# Renumber the running index in column A (0-based) for every data row.
$dates = @("2022-Q1", "2021-Q4", "2021-Q3", "2021-Q2", "2021-Q1", "2020-Q4")
for ($i = 0; $i -lt $dates.Count; $i++) {
    $totals.Range("A" + (2 + $i)).Value = $i
}
